$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '37.175.07'
$c.ClearFormats()
$ws.Range('E2').Value = '  -0.03%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.054.80'
$c.ClearFormats()
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '248.63'
$c.ClearFormats()
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('E6').Value = '  -2.22%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '58.27'
$c.ClearFormats()
$ws.Range('E7').Value = '  -6.42%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -2.51%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0780'
$c.ClearFormats()
$ws.Range('E10').Value = '  -2.98%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.108'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.39%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '15.91'
$c.ClearFormats()
$ws.Range('E12').Value = '  -2.96%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.355.10'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.98%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.836'
$c.ClearFormats()
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('E15').Value = '  +2.36%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.055.96'
$c.ClearFormats()
$ws.Range('E16').Value = '  -1.02%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '18.14'
$c.ClearFormats()
$ws.Range('E17').Value = '  +15.65%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '37.242.51'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.20%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '74.72'
$c.ClearFormats()
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('E20').Value = '  -3.50%  '
$ws.Range('E21').Value = '  -2.66%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '237.43'
$c.ClearFormats()
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('E23').Value = '  -0.10%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.ClearFormats()
$ws.Range('E24').Value = '  +1.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.19'
$c.ClearFormats()
$ws.Range('E25').Value = '  -7.86%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '169.47'
$c.ClearFormats()
$ws.Range('E26').Value = '  -0.63%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.43'
$c.ClearFormats()
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '20.07'
$c.ClearFormats()
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('E29').Value = '  -2.00%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.13'
$c.ClearFormats()
$ws.Range('E30').Value = '  +0.79%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '4.79'
$c.ClearFormats()
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('E32').Value = '  -3.80%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.50'
$c.ClearFormats()
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  -1.50%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  -2.49%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.ClearFormats()
$ws.Range('E39').Value = '  +14.34%  '
$ws.Range('E40').Value = '  +16.57%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.ClearFormats()
$ws.Range('E41').Value = '  -12.64%  '
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('E43').Value = '  -2.49%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '17.20'
$c.ClearFormats()
$ws.Range('E44').Value = '  -5.97%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '96.01'
$c.ClearFormats()
$ws.Range('E45').Value = '  -3.52%  '
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('E47').Value = '  -1.32%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.276.00'
$c.ClearFormats()
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('E49').Value = '  -2.40%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.241.76'
$c.ClearFormats()
$ws.Range('E50').Value = '  -0.81%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '43.71'
$c.ClearFormats()
$ws.Range('E51').Value = '  -1.46%  '
